# aact-784: Refactor so that we will be able to retrieve multiple sheets
# from the Data Definitions spreadsheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Rename the "Data Dictionary" sheet to "Data Definitions"
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Data Definitions"

# ---------------------------------------------------------------------
# 2. Move the selection on the (now renamed) first sheet to B23
# ---------------------------------------------------------------------
$ws1.Activate()
$ws1.Range("B23").Select()

# ---------------------------------------------------------------------
# 3. "General Info" sheet: set explicit column widths for columns
#    A (1), F (6) and G (7).
#    ColumnWidth is specified in "character" units; the values below
#    are chosen so the saved OOXML <col width="..."> lands on the
#    target widths (95.5, 25.5, ~24.33).
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)
$ws2.Columns.Item(1).ColumnWidth = 94.66666666666667
$ws2.Columns.Item(6).ColumnWidth = 24.666666666666668
$ws2.Columns.Item(7).ColumnWidth = 23.498697916666668

# ---------------------------------------------------------------------
# 4. "General Info" sheet: page setup -> landscape, 75% scale
# ---------------------------------------------------------------------
$ws2.PageSetup.Orientation = 2   # xlLandscape
$ws2.PageSetup.Zoom = 75

# ---------------------------------------------------------------------
# 5. Workbook window position (best effort)
# ---------------------------------------------------------------------
$win = $wb.Windows.Item(1)
$win.Left = 620
$win.Top = 3120

Write-Output "edits applied"
